$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# Collapse " con cardinalità " + "0.." + "* da parte dell'utente ... utente
# generico." (which were split across three runs around two <w:proofErr/>
# markers) into a single run/sentence. The visible text does not change,
# only the run/proofErr structure is simplified, so search == replace.
$apos = [char]0x2019
$search1 = " con cardinalità 0..* da parte dell" + $apos + "utente e 1 da parte della prenotazione: un utente generico può effettuare un numero qualsiasi di prenotazioni, mentre ogni prenotazione deve appartenere a un solo utente generico."
$null = $d.Content.Find.Execute($search1, $true, $false, $false, $false, $false, $true, 1, $false, $search1, 2)

# --- Edit 2 -----------------------------------------------------------
# "... esiste solo in presenza di una relazione tra **UtenteGenerico** e
# **Prenotazione**, ..." becomes "... esiste solo in presenza di una
# **Prenotazione**, ...": drop the "relazione tra UtenteGenerico e " text
# (and the bold "UtenteGenerico" run along with it), keeping "Prenotazione"
# bold and intact.
$search2  = " esiste solo in presenza di una relazione tra UtenteGenerico e "
$replace2 = " esiste solo in presenza di una "
$null = $d.Content.Find.Execute($search2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2)

# --- Edit 3 -----------------------------------------------------------
# "..., poiché un utente generico diventa passeggero solo dopo aver
# effettuato una prenotazione, cioè, acquistando un biglietto." becomes
# "..., poiché un utente generico diventa passeggero solo dopo aver
# effettuato una prenotazione, ovvero dopo aver acquistato un biglietto. "
# (note the trailing space), merging the three runs/"cioè," run into one.
$search3  = ", poiché un utente generico diventa passeggero solo dopo aver effettuato una prenotazione, cioè, acquistando un biglietto."
$replace3 = ", poiché un utente generico diventa passeggero solo dopo aver effettuato una prenotazione, ovvero dopo aver acquistato un biglietto. "
$null = $d.Content.Find.Execute($search3, $true, $false, $false, $false, $false, $true, 1, $false, $replace3, 2)
